# Update the "Förändrad" (Changed) date column (C) from 45205 to 45206
# for all data rows (2 through 78) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 78
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
